$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.355.01"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.856.66"
$ws.Range("E3").Value = "  +4.25%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'363.27"
$ws.Range("E5").Value = "  +9.74%  "
$ws.Range("D6").Value = "'116.42"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.554"
$ws.Range("E7").Value = "  +4.17%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  +6.36%  "
$ws.Range("D10").Value = "'42.77"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "'0.0866"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").Value = "'20.19"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'7.90"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("D15").Value = "3.307.99"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "2.858.11"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").Value = "52.487.08"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'14.00"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  +8.44%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'3.18"
$ws.Range("E21").Value = "  +6.30%  "
$ws.Range("D22").Value = "0.0₃0992"
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'70.67"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'272.13"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  +8.21%  "
$ws.Range("D26").Value = "'27.23"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'10.35"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "'0.141"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "'34.64"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'51.24"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.88"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.0447"
$ws.Range("E34").Value = "  +29.76%  "
$ws.Range("D35").Value = "'0.0837"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'5.04"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'18.82"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = "  +10.43%  "
$ws.Range("D42").Value = "'23.68"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").Value = "'126.79"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "'3.40"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "2.081.02"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("D49").Value = "'0.961"
$ws.Range("E49").Value = "  +10.78%  "
$ws.Range("D50").Value = "'5.65"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").Value = "  +0.90%  "
